$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells we are about to rewrite to stay as plain text,
# matching the source file (inline text strings) instead of letting Excel
# auto-detect them as numbers (scientific notation, loss of trailing zeros, etc.)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated values scraped by the GitHub Actions run
$ws.Range("D2").Value = "29.529.84"
$ws.Range("D3").Value = "1.851.77"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("D4").Value = "0.9989"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "241.78"
$ws.Range("E5").Value = "  -1.14%  "
$ws.Range("D6").Value = "0.6282"
$ws.Range("E6").Value = "  -2.17%  "
$ws.Range("D7").Value = "0.9998"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D9").Value = "0.07571"
$ws.Range("D10").Value = "0.2980"
$ws.Range("E10").Value = "  -0.07%  "
$ws.Range("D11").Value = "24.32"
$ws.Range("E11").Value = "  -0.48%  "
$ws.Range("D12").Value = "0.07678"
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").Value = "1.898.49"
$ws.Range("E13").Value = "  +1.45%  "
$ws.Range("D14").Value = "5.017"
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("D15").Value = "0.6858"
$ws.Range("E15").Value = "  -0.88%  "
$ws.Range("D16").Value = "83.91"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "0.000009782"
$ws.Range("E17").Value = "  -0.77%  "
$ws.Range("D18").Value = "2.139.90"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("D19").Value = "6.214"
$ws.Range("E19").Value = "  +1.89%  "
$ws.Range("D20").Value = "29.568.47"
$ws.Range("E20").Value = "  -0.65%  "
$ws.Range("D21").Value = "235.03"
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("E22").Value = "  -1.43%  "
$ws.Range("D23").Value = "0.9999"
$ws.Range("D24").Value = "7.616"
$ws.Range("E24").Value = "  +1.37%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").Value = "155.78"
$ws.Range("E26").Value = "  -2.03%  "
$ws.Range("D27").Value = "0.1391"
$ws.Range("E27").Value = "  -2.18%  "
$ws.Range("D28").Value = "8.428"
$ws.Range("E28").Value = "  -1.14%  "
$ws.Range("D29").Value = "17.73"
$ws.Range("E29").Value = "  -0.97%  "
$ws.Range("D30").Value = "1.480"
$ws.Range("E30").Value = "  -0.74%  "
$ws.Range("E31").Value = "  -6.08%  "
$ws.Range("E32").Value = "  -2.02%  "
$ws.Range("D33").Value = "4.107"
$ws.Range("E33").Value = "  -1.31%  "
$ws.Range("D34").Value = "4.039"
$ws.Range("E34").Value = "  -1.49%  "
$ws.Range("D35").Value = "1.894"
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").Value = "1.173"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("D37").Value = "0.7184"
$ws.Range("E37").Value = "  -1.51%  "
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("D39").Value = "2.801"
$ws.Range("E39").Value = "  -0.81%  "
$ws.Range("D40").Value = "1.237.69"
$ws.Range("E40").Value = "  +3.07%  "
$ws.Range("D41").Value = "0.01779"
$ws.Range("E41").Value = "  -0.55%  "
$ws.Range("D42").Value = "0.9127"
$ws.Range("E42").Value = "  -1.31%  "
$ws.Range("D43").Value = "6.140"
$ws.Range("E43").Value = "  -1.72%  "
$ws.Range("D44").Value = "2.043.35"
$ws.Range("E44").Value = "  +0.51%  "
$ws.Range("D45").Value = "0.9994"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").Value = "101.97"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "67.50"
$ws.Range("E47").Value = "  +1.47%  "
$ws.Range("D48").Value = "7.298"
$ws.Range("E48").Value = "  +9.24%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "9.162"
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").Value = "0.4033"
$ws.Range("E50").Value = "  -0.71%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.00000000117"
$ws.Range("E51").Value = "  -1.18%  "
